{"js": "// Clarify the \"Workaround to copy files into ADLS v2\" paragraph (Page 1,\n// Task 1 - \"Copy source data files to your Azure storage location\").\n//\n// The original paragraph was a single sentence. The update:\n//   1) appends a parenthetical NOTE that the workaround does not work for\n//      hierarchy-enabled v2 storage,\n//   2) inserts a line break after the lead-in sentence, and\n//   3) expands the remaining instructions with guidance about copy time /\n//      verifying file sizes before the \"Then copy...\" sentence.\n\nconst body = context.document.body;\n\nconst oldText =\n  \"Workaround to copy files into ADLS v2: Create v1-enabled storage \" +\n  \"container and copy the source files into the v1 container.  Then copy \" +\n  \"from the v1 storage container into a v2 storage container.  This can \" +\n  \"be done in Storage Explorer copy/paste. \";\n\n// \\v (vertical tab / Word's \"manual line break\" character) becomes a\n// <w:br/> when written back out through Office.js.\nconst newText =\n  \"Workaround to copy files into ADLS v2 (NOTE: This method does not \" +\n  \"work for hierarchy-enabled v2 storage): \\v\" +\n  \"Create v1-enabled storage container and copy the source files into \" +\n  \"the v1 container.  This can be done in Storage Explorer copy/paste.  \" +\n  \"Keep in mind, it will take some time for the files to copy over; \" +\n  \"this does not happen immediately as it appears (folders/files will \" +\n  \"appear, but double-check file sizes to verify the copy completed).  \" +\n  \"Then copy from the v1 storage container into a v2 storage container.  \";\n\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the 'Workaround to copy files into ADLS v2' paragraph text.\");\n}\n\nresults.items[0].insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Clarify the \"Workaround to copy files into ADLS v2\" paragraph (Page 1,\n# Task 1 - \"Copy source data files to your Azure storage location\").\n#\n# The original paragraph was a single sentence. The update:\n#   1) appends a parenthetical NOTE that the workaround does not work for\n#      hierarchy-enabled v2 storage,\n#   2) inserts a line break after the lead-in sentence, and\n#   3) expands the remaining instructions with guidance about copy time /\n#      verifying file sizes before the \"Then copy...\" sentence.\n\n$d = $word.ActiveDocument\n\n$oldText = \"Workaround to copy files into ADLS v2: Create v1-enabled storage container and copy the source files into the v1 container.  Then copy from the v1 storage container into a v2 storage container.  This can be done in Storage Explorer copy/paste. \"\n\n# Chr(11) is Word's manual line-break character (renders as <w:br/>).\n$newText = \"Workaround to copy files into ADLS v2 (NOTE: This method does not work for hierarchy-enabled v2 storage): \" + [char]11 + \"Create v1-enabled storage container and copy the source files into the v1 container.  This can be done in Storage Explorer copy/paste.  Keep in mind, it will take some time for the files to copy over; this does not happen immediately as it appears (folders/files will appear, but double-check file sizes to verify the copy completed).  Then copy from the v1 storage container into a v2 storage container.  \"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
